# Apply crypto price/volume updates from the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.283.74"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "3.130.74"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.126.52"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.41"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "3.646.84"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.118"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "64.198.72"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "3.097.75"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  +7.58%  "
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "0.0₃0770"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.43"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "446.60"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "2.860.24"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.98"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.76"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.08%  "
